# Upper respect and operators doubles
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Salidas Lexer")

# 1. Clear the autofilter criteria (this un-hides the filtered-out rows)
if ($ws.AutoFilterMode) {
    $ws.ShowAllData()
}

# 2. Append new test-case rows (12.out) at rows 199-205
$data = @(
    "<token_cor_izq,1,1>",
    "<id,asd,1,3>",
    "<token_mayor_igual,1,6>",
    "<token_real,87678.9,1,8>",
    "<token_cor_izq,2,1>",
    "<token_cor_der,2,2>",
    "<token_cor_der,2,4>"
)

# Fill columns A, B, C first (row by row) so the shared-string table picks up
# "12.out" and the token strings in source order...
$row = 199
foreach ($d in $data) {
    $ws.Cells.Item($row, 1).Value = "12.out"
    $ws.Cells.Item($row, 2).Value = $d
    $ws.Cells.Item($row, 3).Value = $d
    $row = $row + 1
}

# ...then fill column D last, so ":D" is appended to the shared-string table
# only after all the other new strings.
$row = 199
foreach ($d in $data) {
    $ws.Cells.Item($row, 4).Value = ":D"
    $row = $row + 1
}

# 3. Underline the font of D180 (the autocorrect / upper-respect fix)
$ws.Range("D180").Font.Underline = $true

# 4. Update view: scroll frozen pane and select D180
$ws.Application.ActiveWindow.ScrollRow = 144
$ws.Range("D180").Select()
